# Fix bug: size_primary# = 50 was being incorrectly changed to sizecat_primary# = 3
# because the comparison used ">=" instead of ">".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E4").Value = "size_primary1 >50 & size_primary1 != 888 & size_primary1 != 777"
$ws.Range("E7").Value = "size_primary2 > 50 & size_primary2 != 888 & size_primary2 != 777"

$ws.Range("E8").Select()
